$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10th column) to make room for "Village"
$ws.Columns.Item(10).Insert()

# Match the width of the neighboring columns (G:I) so it merges into the same <col> run
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Set the header value for the new column
$ws.Cells.Item(1, 10).Value = "Village"
